# Update countries & provincias Spain
# Applies the 21-Jun-2020 07:31 -> 08:48 data refresh:
#  - India's daily figures updated
#  - Afganistan overtakes Oman in "Casos totales" (rows swap place / data)
#  - Fiyi / Dominica tie and swap place (values identical, so no visible
#    numeric change, only which country sits in which row)
#  - Islas Turcas y Caicos overtakes Santa Sede
#  - Islas Virgenes Britanicas overtakes Papua Nueva Guinea
#  - El Salvador's daily figures updated
#  - Footer timestamp updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer "last updated" timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 08:48"

# --- Row 7: India (standalone data refresh) -------------------------------
$ws.Range("B7").Value = 411773
$ws.Range("C7").Value = 46
$ws.Range("D7").Value = 228307
$ws.Range("E7").Value = 170185
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 13281

# --- Rows 43/44: Afganistan overtakes Oman --------------------------------
$ws.Range("A43").Value = "Afganistan"
$ws.Range("B43").Value = 28833
$ws.Range("C43").Value = 409
$ws.Range("D43").Value = 8764
$ws.Range("E43").Value = 19488
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 12
$ws.Range("H43").Value = 581

$ws.Range("A44").Value = "Oman"
$ws.Range("B44").Value = 28566
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 14780
$ws.Range("E44").Value = 13658
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 128

# --- Row 85: El Salvador (standalone data refresh) ------------------------
$ws.Range("D85").Value = 2468
$ws.Range("E85").Value = 1909
$ws.Range("G85").Value = 5
$ws.Range("H85").Value = 98

# --- Rows 202/203: Dominica / Fiyi swap places (tied totals) --------------
$ws.Range("A202").Value = "Dominica"
$ws.Range("B202").Value = 18
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 18
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = "Fiyi"
$ws.Range("B203").Value = 18
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# --- Rows 208/209: Islas Turcas y Caicos overtakes Santa Sede -------------
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 11
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("B209").Value = 12
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 12
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# --- Rows 213/214: Islas Virgenes Britanicas overtakes Papua Nueva Guinea -
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
